# Doing Updates for Financials
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AU")

# Earnings Before Interest And Taxes (row 21) - 2011 value now unavailable
$ws.Range("J21").Value = "NA"

# Depreciation (row 83) - 2011 value now unavailable
$ws.Range("J83").Value = "NA"

# Capital Expenditures (row 91) - revised figures across all years
$ws.Range("D91").Value = -829000
$ws.Range("E91").Value = -706000
$ws.Range("F91").Value = -664000
$ws.Range("G91").Value = -844000
$ws.Range("H91").Value = -1363000
$ws.Range("I91").Value = -1925000
$ws.Range("J91").Value = -2485000

# Total Cash Flows From Investing Activities (row 94) - 2011 value now unavailable
$ws.Range("J94").Value = "NA"

# Dividends Paid (row 96) - 2017 value revised
$ws.Range("D96").Value = -39000

# Total Cash Flows From Financing Activities (row 100) - 2011 value now unavailable
$ws.Range("J100").Value = "NA"

# Effect Of Exchange Rate Changes (row 101) - 2011 value now unavailable
$ws.Range("J101").Value = "NA"
